# Update "想去人数" (F column) values across sheets to match the
# newly scraped data output at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 19700
$ws1.Range("F5").Value  = 791
$ws1.Range("F9").Value  = 7456
$ws1.Range("F12").Value = 257
$ws1.Range("F14").Value = 0
$ws1.Range("F20").Value = 0
$ws1.Range("F23").Value = 0
$ws1.Range("F25").Value = 0
$ws1.Range("F26").Value = 315
$ws1.Range("F29").Value = 0
$ws1.Range("F30").Value = 0
$ws1.Range("F34").Value = 2786
$ws1.Range("F35").Value = 24
$ws1.Range("F37").Value = 0
$ws1.Range("F38").Value = 12545
$ws1.Range("F39").Value = 1326
$ws1.Range("F41").Value = 19
$ws1.Range("F44").Value = 0
$ws1.Range("F45").Value = 3980
$ws1.Range("F46").Value = 318
$ws1.Range("F47").Value = 0

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 0

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 19700
$ws4.Range("F5").Value  = 791
$ws4.Range("F8").Value  = 11
$ws4.Range("F9").Value  = 7456
$ws4.Range("F10").Value = 495
$ws4.Range("F12").Value = 257
$ws4.Range("F15").Value = 0
$ws4.Range("F19").Value = 0
$ws4.Range("F20").Value = 386
$ws4.Range("F22").Value = 676
$ws4.Range("F25").Value = 63
$ws4.Range("F26").Value = 315
$ws4.Range("F29").Value = 14
$ws4.Range("F30").Value = 170
$ws4.Range("F32").Value = 0
$ws4.Range("F36").Value = 2786
$ws4.Range("F38").Value = 86
$ws4.Range("F41").Value = 1326
$ws4.Range("F42").Value = 63
$ws4.Range("F43").Value = 19
$ws4.Range("F46").Value = 348
$ws4.Range("F47").Value = 3980
$ws4.Range("F48").Value = 318
